$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = New-Object 'object[,]' 13,5

$rows[0,0]  = "DELL G15 Intel Core i5 13th Gen 13450HX "
$rows[0,1]  = 85990
$rows[0,2]  = 4.2
$rows[0,3]  = "Intel Core i5 13th Gen"
$rows[0,4]  = "1 TB"

$rows[1,0]  = "Acer NITRO LITE 16 Intel Core i7 13th Gen 13620H "
$rows[1,1]  = 69990
$rows[1,2]  = 4.3
$rows[1,3]  = "Intel Core i7 13th Gen"
$rows[1,4]  = "512 GB"

$rows[2,0]  = "Lenovo LOQ 2025 Intel Core i7 14th Gen 14700HX "
$rows[2,1]  = 125990
$rows[2,2]  = 4.2
$rows[2,3]  = "Intel Core i7 14th Gen"
$rows[2,4]  = "1 TB"

$rows[3,0]  = "Lenovo LOQ Essential Intel Core i7 12650HX "
$rows[3,1]  = 79990
$rows[3,2]  = 4.4
$rows[3,3]  = "Intel Core i7"
$rows[3,4]  = "512 GB"

$rows[4,0]  = "Acer NITRO V 16S Intel Core 5 "
$rows[4,1]  = 94990
$rows[4,2]  = 4.8
$rows[4,3]  = "Intel Core 5"
$rows[4,4]  = "512 GB"

$rows[5,0]  = "HP Victus AMD Ryzen 7 Hexa Core 7445HS "
$rows[5,1]  = 63990
$rows[5,2]  = 4.4
$rows[5,3]  = "AMD Ryzen 7 Hexa Core"
$rows[5,4]  = "512 GB"

$rows[6,0]  = "Lenovo LOQ 2025 Intel Core i7 13th Gen 13700HX "
$rows[6,1]  = 118990
$rows[6,2]  = 4.4
$rows[6,3]  = "Intel Core i7 13th Gen"
$rows[6,4]  = "1 TB"

$rows[7,0]  = "HP Victus Intel Core i5 14th Gen 14450HX "
$rows[7,1]  = 89990
$rows[7,2]  = 4.3
$rows[7,3]  = "Intel Core i5 14th Gen"
$rows[7,4]  = "512 GB"

$rows[8,0]  = "HP AMD Ryzen 5 Octa Core 8th Gen "
$rows[8,1]  = 65990
$rows[8,2]  = 3.8
$rows[8,3]  = "AMD Ryzen 5 Octa Core 8th Gen"
$rows[8,4]  = "512 GB"

$rows[9,0]  = "Lenovo LOQ Intel Core i5 12th Gen 12450HX "
$rows[9,1]  = 76990
$rows[9,2]  = 4.5
$rows[9,3]  = "Intel Core i5 12th Gen"
$rows[9,4]  = "512 GB"

$rows[10,0] = "Lenovo IdeaPad Slim 3 AMD Ryzen 7 Octa Core "
$rows[10,1] = 67390
$rows[10,2] = 4.6
$rows[10,3] = "AMD Ryzen 7 Octa Core"
$rows[10,4] = "512 GB"

$rows[11,0] = "ASUS TUF Gaming A16 (2025) AMD Ryzen 7 Hexa Core 7445HS "
$rows[11,1] = 111990
$rows[11,2] = 4.5
$rows[11,3] = "AMD Ryzen 7 Hexa Core"
$rows[11,4] = "1 TB"

$rows[12,0] = "MSI Thin A15 AMD Ryzen 5 Hexa Core 7535HS "
$rows[12,1] = 59990
$rows[12,2] = 4.4
$rows[12,3] = "AMD Ryzen 5 Hexa Core"
$rows[12,4] = "512 GB"

$ws.Range("A2:E14").Value = $rows
